$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = "system, System, backup@backdoor.com"
$ws.Range("G3").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G4").Value2 = "System, backup@backdoor.com"
$ws.Range("G6").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G7").Value2 = "admin@admin.com, System"
$ws.Range("G10").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G11").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G12").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G13").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G14").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G15").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G17").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G18").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G19").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G20").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G21").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G22").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G24").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G26").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G28").Value2 = "system, System, backup@backdoor.com"
$ws.Range("G29").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G30").Value2 = "System, backup@backdoor.com"
$ws.Range("G32").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G33").Value2 = "admin@admin.com, System"
$ws.Range("G36").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G37").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G38").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G39").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G40").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G41").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G43").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G44").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G45").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G46").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G47").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G48").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G50").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G52").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G54").Value2 = "system, System, backup@backdoor.com"
$ws.Range("G55").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G56").Value2 = "System, backup@backdoor.com"
$ws.Range("G58").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G59").Value2 = "admin@admin.com, System"
$ws.Range("G62").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G63").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G64").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G65").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G66").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G67").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G69").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G70").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G71").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G72").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G73").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G74").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G76").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G78").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G83").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G84").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G85").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G86").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G87").Value2 = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G90").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G92").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G93").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G94").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G96").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G99").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G101").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G109").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G110").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G111").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G112").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G113").Value2 = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G116").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G118").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G119").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G120").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G122").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G125").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G127").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G135").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G136").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G137").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G138").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G139").Value2 = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G142").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G144").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G145").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G146").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G148").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G151").Value2 = "dnasr281@gmail.com, System"
$ws.Range("G153").Value2 = "dnasr281@gmail.com, System"
